$d = $word.ActiveDocument

# Replace "PyCharm, " with "Jupyter" in the skills list, then re-insert
# ", " after the bookmark so the final text reads "...Trello, Jupyter, BlueJ"
$d.Content.Find.Execute("PyCharm, ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Jupyter", 2)

$rng = $d.Content
$rng.Find.Execute("Jupyter")
$rng.Collapse(0)
$rng.InsertAfter(", ")
